# DONALD-2: Addition of the categories for parental employment status (EMPLOY_P)
# on the "Categories" sheet. The new categories are inserted right before the
# existing "yes"/"no" (variable index 9) rows, pushing them down by 7 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categories")

# The existing "EMPLOY_P yes/no" rows currently sit at worksheet rows 30-31.
# Insert 7 fresh rows above them so they move down to rows 37-38, then fill
# rows 30-36 with the new employment-status categories.
$ws.Range("A30:A36").EntireRow.Insert()

$categories = @(
    "employed-full time",
    "employed-part time",
    "housewife/-men",
    "retired",
    "unemployed",
    "student",
    "other"
)

for ($i = 0; $i -lt $categories.Length; $i++) {
    $row = 30 + $i
    $ws.Range("A$row").Value = "EMPLOY_P"
    $ws.Range("B$row").Value = $i + 1
    $ws.Range("C$row").Value = $categories[$i]
    # Newly inserted rows picked up the style of the row above (s="1" on
    # column A, matching the bold/no-wrap "variable name" column style used
    # everywhere else) — but the target state leaves these particular cells
    # unstyled, so strip the inherited formatting back off.
    $ws.Range("A$row").ClearFormats()
}

# Restore the view state left behind by the edit: the whole (now relocated)
# last data row selected, with the window scrolled down so row 17 is at the
# top.
$ws.Range("A37:XFD37").Select()
$ws.Application.ActiveWindow.ScrollRow = 17
